# Roll updated device/app identifiers into the "Global" sheet of the
# AI Cross Demo workbook, and move the last-used selection.
#
# The workbook's active (visible) tab is "Action 1"; we temporarily
# activate "Global" so the selection can be set on it, then restore
# the original active sheet so the saved file keeps the same tab
# selected as before.

$wb = $excel.ActiveWorkbook
$originalActive = $wb.ActiveSheet.Name

$ws = $wb.Worksheets.Item("Global")
$ws.Activate()

# device_id for the first (IOS) device row changed
$ws.Range("D4").Value = "8f8fe2aa0724ef5979a590d0f755ece53275b32f"

# Device Name for the first (IOS) device row changed
$ws.Range("H4").Value = "iPhone 6s"

# Column D (device_id) auto-fit width shrank slightly to match the new value
$ws.Columns.Item(4).ColumnWidth = 36.67

# Last cell selected by the user on this sheet
$ws.Range("H5").Select()

# Restore the workbook's originally active sheet/tab
$wb.Worksheets.Item($originalActive).Activate()
